$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values per diff
$ws.Range("C1").Value = 6
$ws.Range("C2").Value = 5
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = 6
$ws.Range("B4").Value = 8
$ws.Range("A5").Value = 8
$ws.Range("C5").Value = 3
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 2
$ws.Range("A7").Value = 4
$ws.Range("C7").Value = 6

# Update active cell selection to A7
$ws.Range("A7").Select()

# Update workbook window position
$excel.Windows.Item(1).Left = 1440
$excel.Windows.Item(1).Top = 880
